$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.787.23"
$ws.Range("D3").Value = "3.625.09"
$ws.Range("E3").Value = "  +4.95%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "588.46"
$ws.Range("E5").Value = "  +2.01%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "193.50"
$ws.Range("E6").Value = "  +3.00%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.639"
$ws.Range("E7").Value = "  +1.43%  "
$ws.Range("D8").Value = "3.618.39"
$ws.Range("E8").Value = "  +5.04%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.00"
$ws.Range("E9").Value = "  +0.08%  "
$ws.Range("E10").Value = "  +5.46%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.669"
$ws.Range("E11").Value = "  +4.27%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "57.61"
$ws.Range("E12").Value = "  -0.23%  "
$ws.Range("E13").Value = "  +10.85%  "
$ws.Range("E14").Value = "  +4.55%  "
$ws.Range("D15").Value = "4.205.76"
$ws.Range("E15").Value = "  +5.22%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "20.32"
$ws.Range("E16").Value = "  +6.74%  "
$ws.Range("D17").Value = "3.616.68"
$ws.Range("E17").Value = "  +4.96%  "
$ws.Range("D18").Value = "70.761.49"
$ws.Range("E18").Value = "  +5.51%  "
$ws.Range("E19").Value = "  +5.10%  "
$ws.Range("E20").Value = "  +2.66%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.06"
$ws.Range("E21").Value = "  +3.71%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "485.87"
$ws.Range("E22").Value = "  -0.58%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "19.34"
$ws.Range("E23").Value = "  +12.04%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.11"
$ws.Range("E24").Value = "  -8.62%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.46"
$ws.Range("E25").Value = "  +3.25%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "90.28"
$ws.Range("E26").Value = "  +0.86%  "
$ws.Range("E27").Value = "  +4.87%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "11.38"
$ws.Range("E28").Value = "  +4.21%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.47"
$ws.Range("E29").Value = "  +4.84%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.94"
$ws.Range("E30").Value = "  +8.13%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "32.56"
$ws.Range("E31").Value = "  +4.04%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.124"
$ws.Range("E32").Value = "  +10.80%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "12.26"
$ws.Range("E33").Value = "  +3.83%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "66.68"
$ws.Range("E34").Value = "  +2.80%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "611.35"
$ws.Range("E35").Value = "  +0.96%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "40.30"
$ws.Range("E36").Value = "  +8.98%  "
$ws.Range("E37").Value = "  +7.03%  "
$ws.Range("E38").Value = "  +5.94%  "
$ws.Range("E39").Value = "  +1.86%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.999"
$ws.Range("E40").Value = "  -0.05%  "
$ws.Range("E41").Value = "  +16.06%  "
$ws.Range("E42").Value = "  +3.45%  "
$ws.Range("D43").Value = "3.322.98"
$ws.Range("E43").Value = "  +4.12%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.21"
$ws.Range("E44").Value = "  +20.75%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.15"
$ws.Range("E45").Value = "  +9.64%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0456"
$ws.Range("E46").Value = "  +6.52%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.75"
$ws.Range("E47").Value = "  +13.39%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.36"
$ws.Range("E48").Value = "  +4.41%  "
$ws.Range("E49").Value = "  +2.56%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.997"
$ws.Range("E50").Value = "  -0.21%  "
$ws.Range("E51").Value = "  +1.75%  "
